$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as text so numeric-looking price strings
# (e.g. "187.75") are not coerced into real numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "66.061.10"
$ws.Range("E2").Value = "  -2.44%  "
$ws.Range("D3").Value = "3.413.71"
$ws.Range("E3").Value = "  -5.32%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("D5").Value = "187.75"
$ws.Range("E5").Value = "  -5.58%  "
$ws.Range("D6").Value = "534.76"
$ws.Range("E6").Value = "  -4.29%  "
$ws.Range("D7").Value = "0.621"
$ws.Range("E7").Value = "  +1.37%  "
$ws.Range("D8").Value = "3.406.34"
$ws.Range("E8").Value = "  -5.33%  "
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  -0.14%  "
$ws.Range("D10").Value = "0.636"
$ws.Range("E10").Value = "  -5.06%  "
$ws.Range("D11").Value = "58.48"
$ws.Range("E11").Value = "  +0.55%  "
$ws.Range("E12").Value = "  -9.58%  "
$ws.Range("D13").Value = "0.0000259"
$ws.Range("E13").Value = "  -10.21%  "
$ws.Range("D14").Value = "9.46"
$ws.Range("E14").Value = "  -5.18%  "
$ws.Range("D15").Value = "3.937.75"
$ws.Range("E15").Value = "  -5.84%  "
$ws.Range("E16").Value = "  -2.09%  "
$ws.Range("D17").Value = "3.405.33"
$ws.Range("E17").Value = "  -5.53%  "
$ws.Range("D18").Value = "65.649.62"
$ws.Range("E18").Value = "  -2.84%  "
$ws.Range("D19").Value = "17.75"
$ws.Range("E19").Value = "  -6.17%  "
$ws.Range("D20").Value = "11.40"
$ws.Range("E20").Value = "  -7.17%  "
$ws.Range("D21").Value = "0.993"
$ws.Range("E21").Value = "  -7.81%  "
$ws.Range("D22").Value = "385.68"
$ws.Range("E22").Value = "  -3.43%  "
$ws.Range("D23").Value = "83.91"
$ws.Range("E23").Value = "  -1.34%  "
$ws.Range("D24").Value = "3.80"
$ws.Range("E24").Value = "  -7.81%  "
$ws.Range("D25").Value = "11.21"
$ws.Range("E25").Value = "  -14.34%  "
$ws.Range("D26").Value = "3.74"
$ws.Range("E26").Value = "  -2.52%  "
$ws.Range("D27").Value = "2.72"
$ws.Range("E27").Value = "  -7.48%  "
$ws.Range("D28").Value = "11.76"
$ws.Range("E28").Value = "  -5.56%  "
$ws.Range("D29").Value = "8.60"
$ws.Range("E29").Value = "  -7.71%  "
$ws.Range("D30").Value = "682.55"
$ws.Range("E30").Value = "  +1.70%  "
$ws.Range("D31").Value = "30.08"
$ws.Range("E31").Value = "  -4.42%  "
$ws.Range("D32").Value = "6.81"
$ws.Range("E32").Value = "  -20.05%  "
$ws.Range("B33").Value = "Cosmos"
$ws.Range("C33").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D33").Value = "11.35"
$ws.Range("E33").Value = "  -6.73%  "
$ws.Range("B34").Value = "OKB"
$ws.Range("C34").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D34").Value = "62.09"
$ws.Range("E34").Value = "  -2.35%  "
$ws.Range("E35").Value = "  -5.09%  "
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  -0.12%  "
$ws.Range("D37").Value = "37.20"
$ws.Range("E37").Value = "  -12.44%  "
$ws.Range("D38").Value = "0.383"
$ws.Range("E38").Value = "  -11.72%  "
$ws.Range("D39").Value = "0.998"
$ws.Range("E39").Value = "  -0.05%  "
$ws.Range("E40").Value = "  -7.23%  "
$ws.Range("D41").Value = "2.930.18"
$ws.Range("E41").Value = "  -9.69%  "
$ws.Range("D42").Value = "2.82"
$ws.Range("E42").Value = "  -11.43%  "
$ws.Range("B43").Value = "PEPE"
$ws.Range("C43").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D43").Value = "0.0₃0639"
$ws.Range("E43").Value = "  -17.44%  "
$ws.Range("B44").Value = "Fetch.AI"
$ws.Range("C44").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D44").Value = "2.42"
$ws.Range("E44").Value = "  -14.78%  "
$ws.Range("D45").Value = "2.65"
$ws.Range("E45").Value = "  -1.99%  "
$ws.Range("D46").Value = "0.0392"
$ws.Range("E46").Value = "  -5.80%  "
$ws.Range("D47").Value = "0.128"
$ws.Range("E47").Value = "  -2.09%  "
$ws.Range("D48").Value = "2.95"
$ws.Range("E48").Value = "  -4.65%  "
$ws.Range("D49").Value = "133.46"
$ws.Range("E49").Value = "  -3.60%  "
$ws.Range("D50").Value = "2.40"
$ws.Range("E50").Value = "  -21.11%  "
$ws.Range("D51").Value = "2.59"
$ws.Range("E51").Value = "  -5.14%  "

# Restore the default style on column D so no stray number-format
# style id is left attached to the cells.
$ws.Range("D2:D51").Style = "Normal"

